# Minor changes to audibility docs
#
# 1. Bump the cached "fixed" date shown by the Date placeholder on the
#    slide master and every slide layout: 14/04/2015 -> 20/04/2015.
# 2. On slide 2, change the "(approved)" caption to "(final)" on the two
#    Q1/Q2 invoice rectangles.
# 3. Add two presentation-level slide guides (horizontal @ 2160,
#    vertical @ 2880) - attempted defensively in case the host supports
#    it; wrapped so it can never abort the rest of the script.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text on the slide master + all slide layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container, [string]$newText) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

Update-DatePlaceholder $p.SlideMaster "20/04/2015"

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i) "20/04/2015"
}

# ---------------------------------------------------------------------
# 2) "(approved)" -> "(final)" on slide 2.
# ---------------------------------------------------------------------
function Update-SubText($shape, [string]$find, [string]$replace) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($find)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $find.Length)
        $sub.Text = $replace
    }
}

$slide2 = $p.Slides.Item(2)
Update-SubText $slide2.Shapes.Item("Rectangle 38") "(approved)" "(final)"
Update-SubText $slide2.Shapes.Item("Rectangle 39") "(approved)" "(final)"

# ---------------------------------------------------------------------
# 3) Presentation slide guides: horizontal guide @ 2160, vertical @ 2880.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        $hGuide = $guides.Add(1, 2160)
        if ($hGuide -ne $null) { $hGuide.Position = 2160 }
        $vGuide = $guides.Add(2, 2880)
        if ($vGuide -ne $null) { $vGuide.Position = 2880 }
    }
} catch {
    # Guides collection isn't available in this host; nothing more to do.
}
